$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.9380533998416762
$ws.Range("E2").Value = 1.155022846258058

$ws.Range("C3").Value = -4.857241224140929
$ws.Range("E3").Value = -1.691674259276643

$ws.Range("C4").Value = 0.3625742673738941
$ws.Range("E4").Value = -1.432000573345915

$ws.Range("C5").Value = 0.9965309787904442
$ws.Range("E5").Value = -0.02552935725171901

$ws.Range("C6").Value = 0.6836026627130565
$ws.Range("E6").Value = 1.279644059586871

$ws.Range("C7").Value = 0.2336505480021955
$ws.Range("E7").Value = 0.4758549421990166

$ws.Range("C8").Value = -0.1446844164011307
$ws.Range("E8").Value = -0.03923323971219972

$ws.Range("C9").Value = -0.1588690085687849
$ws.Range("E9").Value = -0.4617525814883283

$ws.Range("C10").Value = -0.5438176183081733
$ws.Range("E10").Value = 0.01247916696665019

$ws.Range("C11").Value = -0.006876704825709012
$ws.Range("E11").Value = -0.05033169102144353

$ws.Range("C12").Value = 0.6718983809452572
$ws.Range("E12").Value = 0.08750765859864007

$ws.Range("C13").Value = -0.6347897325981511
$ws.Range("E13").Value = 0.01241557525981651

$ws.Range("C14").Value = -0.7158018152081724
$ws.Range("E14").Value = -0.76158050880345

$ws.Range("C15").Value = 1.22331349480691
$ws.Range("E15").Value = -0.5204428773059266

$ws.Range("C16").Value = -1.508346016334061
$ws.Range("E16").Value = -1.034042971854776

$ws.Range("C17").Value = 0.2857158074419441
$ws.Range("E17").Value = -0.2434890887128005

$ws.Range("C18").Value = 1.050311853611596
$ws.Range("E18").Value = 0.6547560647617745

$ws.Range("C19").Value = -1.467296258526263
$ws.Range("E19").Value = -0.02957675682233596
